$d = $word.ActiveDocument

# --- Paragraph 2 contains the field text: {m:'prefix<CR>suffix'}
# Target layout (per TokenIteratorFieldRewriterSplit):
#   "{"  "m"  ":'prefix"  [bookmark]  "\r"  "suffix"  " "  "'"  "}"
#
# Step 1: split the "{m" run into "{" and "m" by toggling a character
# format on just the "{" character; Word always starts a new run at a
# formatting boundary even when the value is later restored.
$brace = $d.Range(36, 37)
$brace.Font.Bold = $true
$brace.Font.Bold = $false

# Step 2: insert the missing space before the closing quote, turning
# "...suffix'}" into "...suffix '}".
$insertionPoint = $d.Range(54, 54)
$insertionPoint.InsertBefore(" ")

# Step 3: the insertion collapses the neighbouring same-format runs
# together, so re-split "\r" from "suffix", and split the new
# " '}" tail into three separate runs: " ", "'", "}".
$cr = $d.Range(46, 48)
$cr.Font.Bold = $true
$cr.Font.Bold = $false

$space = $d.Range(54, 55)
$space.Font.Bold = $true
$space.Font.Bold = $false

$quote = $d.Range(55, 56)
$quote.Font.Bold = $true
$quote.Font.Bold = $false

Write-Output "done"
